$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 4 formatting touch-up ----------------------------------------
# D4 / E4 / G4 were using a "Noto Sans CJK SC" font style that had a
# duplicate twin in the stylesheet; re-asserting the font collapses them
# onto the shared style.
$ws.Range("D4").Font.Name = "Noto Sans CJK SC"
$ws.Range("E4").Font.Name = "Noto Sans CJK SC"
$ws.Range("G4").Font.Name = "Noto Sans CJK SC"

# F4 keeps its wrapped-text formatting, but also gets the normalised font.
$ws.Range("F4").WrapText = $true
$ws.Range("F4").Font.Name = "Noto Sans CJK SC"

# --- New note in F5 -----------------------------------------------------
# "dof" typed in the Latin font, followed by the Japanese continuation in
# the CJK font, matching the workbook's existing mixed-run convention.
$ws.Range("F5").Value = "dofも初期位置ランダム要素を加えてみたらどうなるだろうか？"
$ws.Range("F5").Characters(1, 3).Font.Name = "Arial"
$ws.Range("F5").Characters(4, 100).Font.Name = "Noto Sans CJK SC"

# Leave the selection where the user's cursor ended up after the edit.
$ws.Range("F6").Select() | Out-Null
